# Updated cryptos list on Thu Jul 27 03:44:04 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as Text so numeric-looking strings
    # (e.g. "0.9996") are not silently converted into real numbers,
    # then restore the default "Normal" style so no stray formatting
    # is left behind on the cell.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "29.436.80"
Set-TextValue "E2" "  +0.74%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.878.30"

# Row 4 - TetherUSD
Set-TextValue "D4" "0.9996"
Set-TextValue "E4" "  +0.07%  "

# Row 5 - XRP
Set-TextValue "D5" "0.7207"
Set-TextValue "E5" "  +1.43%  "

# Row 6 - BNB
Set-TextValue "D6" "240.20"
Set-TextValue "E6" "  +0.84%  "

# Row 7 - USDC
Set-TextValue "D7" "0.9998"
Set-TextValue "E7" "  +0.03%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.07840"
Set-TextValue "E8" "  -2.12%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.3106"
Set-TextValue "E9" "  +2.34%  "

# Row 10 - Solana
Set-TextValue "D10" "24.96"
Set-TextValue "E10" "  +5.85%  "

# Row 11 - TRON
Set-TextValue "D11" "0.08262"
Set-TextValue "E11" "  +0.89%  "

# Row 12 - Polygon
Set-TextValue "D12" "0.7281"
Set-TextValue "E12" "  +3.33%  "

# Row 13 - Polkadot
Set-TextValue "D13" "5.285"
Set-TextValue "E13" "  +2.12%  "

# Row 14 - WrappedEther
Set-TextValue "D14" "1.844.49"
Set-TextValue "E14" "  -0.15%  "

# Row 15 - Litecoin
Set-TextValue "D15" "91.28"
Set-TextValue "E15" "  +1.74%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "29.383.09"
Set-TextValue "E16" "  +0.61%  "

# Row 17 - Uniswap
Set-TextValue "D17" "5.937"
Set-TextValue "E17" "  +1.63%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "245.39"
Set-TextValue "E18" "  +2.97%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.000007908"
Set-TextValue "E19" "  +0.38%  "

# Row 20 - Avalanche
Set-TextValue "D20" "13.32"

# Row 21 - WrappedliquidstakedEther2.0
Set-TextValue "D21" "2.117.37"
Set-TextValue "E21" "  +1.41%  "

# Row 22 - Dai
Set-TextValue "D22" "0.9991"
Set-TextValue "E22" "  +0.14%  "

# Row 23 - Chainlink
Set-TextValue "D23" "7.949"
Set-TextValue "E23" "  +6.81%  "

# Row 24 - BinanceUSD
Set-TextValue "D24" "0.9994"
Set-TextValue "E24" "  +0.00%  "

# Row 25 - Stellar
Set-TextValue "D25" "0.1588"
Set-TextValue "E25" "  +10.06%  "

# Row 26 - Monero
Set-TextValue "D26" "164.01"
Set-TextValue "E26" "  +1.05%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.033"
Set-TextValue "E27" "  +1.11%  "

# Row 28 - EthereumClassic
Set-TextValue "E28" "  +1.45%  "

# Row 29 - Toncoin
Set-TextValue "E29" "  -4.67%  "

# Row 30 / Row 31 - Filecoin and PancakeSwap swap places
Set-TextValue "B30" "PancakeSwap"
Set-TextValue "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "1.485"
Set-TextValue "E30" "  +0.33%  "

Set-TextValue "B31" "Filecoin"
Set-TextValue "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "4.397"
Set-TextValue "E31" "  +0.54%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "4.154"
Set-TextValue "E32" "  +3.51%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.05286"
Set-TextValue "E33" "  +1.41%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.946"
Set-TextValue "E34" "  +0.74%  "

# Row 35 - ARBITRUM
Set-TextValue "E35" "  +3.37%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.7238"

# Row 38 - VeChain
Set-TextValue "D38" "0.01869"
Set-TextValue "E38" "  +0.57%  "

# Row 39 - Maker
Set-TextValue "D39" "1.241.51"
Set-TextValue "E39" "  +9.80%  "

# Row 40 - MXToken
Set-TextValue "D40" "2.719"
Set-TextValue "E40" "  -0.23%  "

# Row 41 - TrustWalletToken
Set-TextValue "D41" "0.9090"
Set-TextValue "E41" "  -2.22%  "

# Row 42 - Aave
Set-TextValue "D42" "73.11"
Set-TextValue "E42" "  +4.23%  "

# Row 43 - FraxShare
Set-TextValue "D43" "6.090"
Set-TextValue "E43" "  +3.96%  "

# Row 44 - PaxDollar
Set-TextValue "D44" "0.9999"

# Row 45 - Quant
Set-TextValue "D45" "103.76"
Set-TextValue "E45" "  +0.73%  "

# Row 46 - Mantle
Set-TextValue "D46" "0.5331"
Set-TextValue "E46" "  -0.17%  "

# Rows 47/48/49 - SynthetixNetwork, BabyDogeCoin, RenderToken rotate
Set-TextValue "B47" "BabyDogeCoin"
Set-TextValue "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.00000000121"
Set-TextValue "E47" "  +1.72%  "

Set-TextValue "B48" "RenderToken"
Set-TextValue "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.766"
Set-TextValue "E48" "  +0.04%  "

Set-TextValue "B49" "SynthetixNetwork"
Set-TextValue "C49" "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue "D49" "2.924"
Set-TextValue "E49" "  +12.83%  "

# Row 50 - TheSandbox
Set-TextValue "D50" "0.4338"
Set-TextValue "E50" "  +1.63%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "9.300"
Set-TextValue "E51" "  +1.52%  "
